$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Ping%Zhang%NULL%1,                Zhigang%He%NULL%1,                Gang%Yu%NULL%1,                Dan%Peng%NULL%1,                Yikuan%Feng%NULL%1,                Jianmin%Ling%NULL%1,                Ye%Wang%NULL%1,                Shusheng%Li%NULL%0,                Yi%Bian%NULL%1]"
$ws.Range("E3").Value = "[Tao%Li%NULL%1,                Yalan%Zhang%NULL%2,                Yalan%Zhang%NULL%0,                Cheng%Gong%NULL%1,                Jing%Wang%NULL%0,                Bao%Liu%NULL%1,                Li%Shi%NULL%1,                Jun%Duan%junjununun@163.com%1]"
$ws.Range("E5").Value = "[Gaoli%Liu%NULL%1,                Shaowen%Zhang%NULL%2,                Shaowen%Zhang%NULL%0,                Zhangfan%Mao%NULL%1,                Weixing%Wang%13392186@qq.com%1,                Haifeng%Hu%NULL%1]"

$wb.Save()
